# Add library preparer information to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B = libraryPreparer, Column E = purpose.
# Replace the old "Retrofitted_2288" placeholder values (rows 2-25) with the
# updated purpose for this run and the real preparer initials.
$ws.Range("E2:E25").Value = "fullRNASEQ"
$ws.Range("B2:B25").Value = "H.BROWN"

# Reflect the selection/scroll state left behind after filling column B.
$ws.Range("A7").Select()
$ws.Range("B3:B25").Select()
